# Auto-generated edit script: refreshes cached market-price / profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to
# match a scheduled data-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3956.1428
$ws.Range("I69").Value = 3942.6
$ws.Range("J69").Value = 3990
$ws.Range("K69").Value = 11827.8
$ws.Range("L69").Value = 11970
$ws.Range("M69").Value = -10953.8
$ws.Range("N69").Value = -13718
$ws.Range("H72").Value = 3956.1428
$ws.Range("I72").Value = 3942.6
$ws.Range("J72").Value = 3990
$ws.Range("K72").Value = 35483.4
$ws.Range("L72").Value = 35910
$ws.Range("M72").Value = -31115.4
$ws.Range("N72").Value = -44646
$ws.Range("H82").Value = 560.5
$ws.Range("I82").Value = 560.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1681.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1275.5
$ws.Range("H85").Value = 560.5
$ws.Range("I85").Value = 560.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1681.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -277.5
$ws.Range("H113").Value = 3249298.5
$ws.Range("I113").Value = 8405461
$ws.Range("J113").Value = 2825.963
$ws.Range("K113").Value = 8405461
$ws.Range("L113").Value = 2825.963
$ws.Range("M113").Value = -8402207
$ws.Range("N113").Value = -9333.963
$ws.Range("H132").Value = 24595158
$ws.Range("I132").Value = 27278194
$ws.Range("J132").Value = 667.1667
$ws.Range("K132").Value = 81834582
$ws.Range("L132").Value = 2001.5001
$ws.Range("M132").Value = -81832052
$ws.Range("N132").Value = -7061.5001
$ws.Range("H137").Value = 2579.638
$ws.Range("I137").Value = 2405.1428
$ws.Range("J137").Value = 3037.6875
$ws.Range("K137").Value = 7215.428400000001
$ws.Range("L137").Value = 9113.0625
$ws.Range("M137").Value = -4665.428400000001
$ws.Range("N137").Value = -14213.0625
$ws.Range("H141").Value = 3403.5454
$ws.Range("I141").Value = 1691.75
$ws.Range("K141").Value = 5075.25
$ws.Range("M141").Value = 104.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1260.4
$ws.Range("I2").Value = 985.1429000000001
$ws.Range("J2").Value = 1902.6666
$ws.Range("K2").Value = 985.1429000000001
$ws.Range("L2").Value = 1902.6666
$ws.Range("M2").Value = -872.1429000000001
$ws.Range("N2").Value = -2128.6666
$ws.Range("H32").Value = 5381654.5
$ws.Range("I32").Value = 4309.1865
$ws.Range("J32").Value = 111136110
$ws.Range("K32").Value = 4309.1865
$ws.Range("L32").Value = 111136110
$ws.Range("M32").Value = -4022.1865
$ws.Range("N32").Value = -111136684
$ws.Range("H45").Value = 15921785
$ws.Range("I45").Value = 27862004
$ws.Range("K45").Value = 27862004
$ws.Range("M45").Value = -27861627
$ws.Range("H61").Value = 3290722.8
$ws.Range("I61").Value = 4465094
$ws.Range("K61").Value = 4465094
$ws.Range("M61").Value = -4464882
$ws.Range("H116").Value = 1260.4
$ws.Range("I116").Value = 985.1429000000001
$ws.Range("J116").Value = 1902.6666
$ws.Range("K116").Value = 985.1429000000001
$ws.Range("L116").Value = 1902.6666
$ws.Range("M116").Value = 1308.8571
$ws.Range("N116").Value = -6490.6666
$ws.Range("H122").Value = 1823.8096
$ws.Range("I122").Value = 1618.75
$ws.Range("J122").Value = 2480
$ws.Range("K122").Value = 4856.25
$ws.Range("L122").Value = 7440
$ws.Range("M122").Value = -2406.25
$ws.Range("N122").Value = -12340
$ws.Range("H135").Value = 60475.5
$ws.Range("J135").Value = 60475.5
$ws.Range("L135").Value = 60475.5
$ws.Range("N135").Value = -70615.5
$ws.Range("H136").Value = 3290722.8
$ws.Range("I136").Value = 4465094
$ws.Range("K136").Value = 13395282
$ws.Range("M136").Value = -13392732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1260.4
$ws.Range("I3").Value = 985.1429000000001
$ws.Range("J3").Value = 1902.6666
$ws.Range("K3").Value = 985.1429000000001
$ws.Range("L3").Value = 1902.6666
$ws.Range("M3").Value = -871.1429000000001
$ws.Range("N3").Value = -2130.6666
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H94").Value = 821.8
$ws.Range("I94").Value = 821.8
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 821.8
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -370.8
$ws.Range("H99").Value = 869.1667
$ws.Range("I99").Value = 792.2222
$ws.Range("J99").Value = 1100
$ws.Range("K99").Value = 792.2222
$ws.Range("L99").Value = 1100
$ws.Range("M99").Value = 705.7778
$ws.Range("N99").Value = -4096
$ws.Range("H107").Value = 27779756
$ws.Range("I107").Value = 45455710
$ws.Range("J107").Value = 3257.1428
$ws.Range("K107").Value = 45455710
$ws.Range("L107").Value = 3257.1428
$ws.Range("M107").Value = -45453790
$ws.Range("N107").Value = -7097.1428
$ws.Range("H109").Value = 30669.6
$ws.Range("J109").Value = 30669.6
$ws.Range("L109").Value = 30669.6
$ws.Range("N109").Value = -33443.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2671.2856
$ws.Range("I16").Value = 1800
$ws.Range("J16").Value = 2816.5
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 2816.5
$ws.Range("M16").Value = -1513
$ws.Range("N16").Value = -3390.5
$ws.Range("H31").Value = 1225.841
$ws.Range("I31").Value = 728.1429000000001
$ws.Range("J31").Value = 2096.8125
$ws.Range("K31").Value = 728.1429000000001
$ws.Range("L31").Value = 2096.8125
$ws.Range("M31").Value = -433.1429000000001
$ws.Range("N31").Value = -2686.8125
$ws.Range("H34").Value = 1225.841
$ws.Range("I34").Value = 728.1429000000001
$ws.Range("J34").Value = 2096.8125
$ws.Range("K34").Value = 728.1429000000001
$ws.Range("L34").Value = 2096.8125
$ws.Range("M34").Value = -526.1429000000001
$ws.Range("N34").Value = -2500.8125
$ws.Range("H113").Value = 2671.2856
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 2816.5
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2816.5
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7156.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 126
$ws.Range("J23").Value = 135.2
$ws.Range("L23").Value = 405.6
$ws.Range("N23").Value = -875.5999999999999
$ws.Range("H129").Value = 1293.5652
$ws.Range("I129").Value = 755.63635
$ws.Range("J129").Value = 1786.6666
$ws.Range("K129").Value = 2266.90905
$ws.Range("L129").Value = 5359.9998
$ws.Range("M129").Value = 2733.09095
$ws.Range("N129").Value = -15359.9998
$ws.Range("H131").Value = 2829.34
$ws.Range("I131").Value = 3537.3157
$ws.Range("J131").Value = 2663.2715
$ws.Range("K131").Value = 10611.9471
$ws.Range("L131").Value = 7989.814499999999
$ws.Range("M131").Value = -5571.947100000001
$ws.Range("N131").Value = -18069.8145
$ws.Range("H139").Value = 2093.5715
$ws.Range("J139").Value = 5000
$ws.Range("L139").Value = 15000
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 25980
$ws.Range("J21").Value = 25980
$ws.Range("L21").Value = 25980
$ws.Range("N21").Value = -26326
$ws.Range("H30").Value = 25980
$ws.Range("J30").Value = 25980
$ws.Range("L30").Value = 25980
$ws.Range("N30").Value = -26190
$ws.Range("H57").Value = 14663.765
$ws.Range("J57").Value = 14663.765
$ws.Range("L57").Value = 14663.765
$ws.Range("N57").Value = -16303.765
$ws.Range("H80").Value = 1983
$ws.Range("I80").Value = 1850
$ws.Range("K80").Value = 1850
$ws.Range("M80").Value = -852
$ws.Range("H83").Value = 1983
$ws.Range("I83").Value = 1850
$ws.Range("K83").Value = 9250
$ws.Range("M83").Value = -4258
$ws.Range("H113").Value = 1744.4286
$ws.Range("I113").Value = 1535.1666
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1535.1666
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 634.8334
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 62502130
$ws.Range("I40").Value = 2834.6667
$ws.Range("K40").Value = 2834.6667
$ws.Range("M40").Value = -2698.6667
$ws.Range("H61").Value = 2032.5714
$ws.Range("I61").Value = 1155.75
$ws.Range("J61").Value = 3201.6667
$ws.Range("K61").Value = 1155.75
$ws.Range("L61").Value = 3201.6667
$ws.Range("M61").Value = -953.75
$ws.Range("N61").Value = -3605.6667
$ws.Range("H113").Value = 2032.5714
$ws.Range("I113").Value = 1155.75
$ws.Range("J113").Value = 3201.6667
$ws.Range("K113").Value = 1155.75
$ws.Range("L113").Value = 3201.6667
$ws.Range("M113").Value = 1014.25
$ws.Range("N113").Value = -7541.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11058.25
$ws.Range("J14").Value = 11058.25
$ws.Range("L14").Value = 11058.25
$ws.Range("N14").Value = -11394.25
$ws.Range("H45").Value = 17813
$ws.Range("J45").Value = 17813
$ws.Range("L45").Value = 17813
$ws.Range("N45").Value = -18795
$ws.Range("H52").Value = 3500
$ws.Range("I52").Value = 500
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 500
$ws.Range("L52").Value = 5000
$ws.Range("M52").Value = -274
$ws.Range("N52").Value = -5452

# Rows whose HQ price dropped to/rose from 0 gain or lose their N column
# entirely in the source export, so clear it explicitly here.
$wb.Worksheets.Item("ALC").Range("N82").ClearContents()
$wb.Worksheets.Item("ALC").Range("N85").ClearContents()
$wb.Worksheets.Item("BSM").Range("N94").ClearContents()
$wb.Worksheets.Item("BSM").Range("N32").ClearContents()
